$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.1309698538606213
$ws.Range("J2").Value = 0.1843806287874228
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.119639
$ws.Range("N2").Value = 0.358917
$ws.Range("O2").Value = 0.02933906125009379
$ws.Range("P2").Value = 0.03107835843382883
$ws.Range("Q2").Value = 0.03836519644533333
$ws.Range("R2").Value = 0.345286768008
$ws.Range("S2").Value = 0.003842532564332602
$ws.Range("T2").Value = 0.005730247269710265

# Row 3
$ws.Range("I3").Value = 0.1309698538606213
$ws.Range("J3").Value = 0.1843806287874228
$ws.Range("O3").Value = 0.05256748359289284
$ws.Range("P3").Value = 0.05568382311683939
$ws.Range("S3").Value = 0.006884755643981784
$ws.Range("T3").Value = 0.01026701831957048

# Row 4
$ws.Range("I4").Value = 0.1309698538606213
$ws.Range("J4").Value = 0.1843806287874228
$ws.Range("M4").Value = 1.018537666666667
$ws.Range("N4").Value = 3.055613
$ws.Range("O4").Value = 0.2497759007335481
$ws.Range("P4").Value = 0.2645832770503125
$ws.Range("Q4").Value = 0.3266192267457778
$ws.Range("R4").Value = 2.939573040712
$ws.Range("S4").Value = 0.03271311321697785
$ws.Range("T4").Value = 0.04878403098917353

# Row 5
$ws.Range("I5").Value = 0.1309698538606213
$ws.Range("J5").Value = 0.1843806287874228
$ws.Range("M5").Value = 0.684642
$ws.Range("N5").Value = 1.369284
$ws.Range("O5").Value = 0.1678946963146358
$ws.Range("P5").Value = 0.1185652921140734
$ws.Range("Q5").Value = 0.219547345136
$ws.Range("R5").Value = 1.317284070816
$ws.Range("S5").Value = 0.02198914384030125
$ws.Range("T5").Value = 0.02186114311235732

# Row 6
$ws.Range("I6").Value = 0.1309698538606213
$ws.Range("J6").Value = 0.1843806287874228
$ws.Range("M6").Value = 2.040627333333334
$ws.Range("N6").Value = 6.121882
$ws.Range("O6").Value = 0.5004228581088294
$ws.Range("P6").Value = 0.5300892492849458
$ws.Range("Q6").Value = 0.6543774899075556
$ws.Range("R6").Value = 5.889397409168001
$ws.Range("S6").Value = 0.06554030859502782
$ws.Range("T6").Value = 0.09773818909661125

# Row 7
$ws.Range("G7").Value = 2.127787
$ws.Range("H7").Value = 4.255574
$ws.Range("I7").Value = 0.8690301461393787
$ws.Range("J7").Value = 0.8156193712125771
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.119639
$ws.Range("N7").Value = 0.358917
$ws.Range("O7").Value = 0.02933906125009379
$ws.Range("P7").Value = 0.03107835843382883
$ws.Range("Q7").Value = 0.254566308893
$ws.Range("R7").Value = 1.527397853358
$ws.Range("S7").Value = 0.02549652868576119
$ws.Range("T7").Value = 0.02534811116411856

# Row 8
$ws.Range("G8").Value = 2.127787
$ws.Range("H8").Value = 4.255574
$ws.Range("I8").Value = 0.8690301461393787
$ws.Range("J8").Value = 0.8156193712125771
$ws.Range("O8").Value = 0.05256748359289284
$ws.Range("P8").Value = 0.05568382311683939
$ws.Range("Q8").Value = 0.45611242132
$ws.Range("R8").Value = 2.73667452792
$ws.Range("S8").Value = 0.04568272794891105
$ws.Range("T8").Value = 0.04541680479726891

# Row 9
$ws.Range("G9").Value = 2.127787
$ws.Range("H9").Value = 4.255574
$ws.Range("I9").Value = 0.8690301461393787
$ws.Range("J9").Value = 0.8156193712125771
$ws.Range("M9").Value = 1.018537666666667
$ws.Range("N9").Value = 3.055613
$ws.Range("O9").Value = 0.2497759007335481
$ws.Range("P9").Value = 0.2645832770503125
$ws.Range("Q9").Value = 2.167231206143667
$ws.Range("R9").Value = 13.003387236862
$ws.Range("S9").Value = 0.2170627875165702
$ws.Range("T9").Value = 0.215799246061139

# Row 10
$ws.Range("G10").Value = 2.127787
$ws.Range("H10").Value = 4.255574
$ws.Range("I10").Value = 0.8690301461393787
$ws.Range("J10").Value = 0.8156193712125771
$ws.Range("M10").Value = 0.684642
$ws.Range("N10").Value = 1.369284
$ws.Range("O10").Value = 0.1678946963146358
$ws.Range("P10").Value = 0.1185652921140734
$ws.Range("Q10").Value = 1.456772347254
$ws.Range("R10").Value = 5.827089389016
$ws.Range("S10").Value = 0.1459055524743346
$ws.Range("T10").Value = 0.09670414900171606

# Row 11
$ws.Range("G11").Value = 2.127787
$ws.Range("H11").Value = 4.255574
$ws.Range("I11").Value = 0.8690301461393787
$ws.Range("J11").Value = 0.8156193712125771
$ws.Range("M11").Value = 2.040627333333334
$ws.Range("N11").Value = 6.121882
$ws.Range("O11").Value = 0.5004228581088294
$ws.Range("P11").Value = 0.5300892492849458
$ws.Range("Q11").Value = 4.342020311711334
$ws.Range("R11").Value = 26.052121870268
$ws.Range("S11").Value = 0.4348825495138016
$ws.Range("T11").Value = 0.4323510601883346

